$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Fitness" column (C) values for rows 2-12 as per the fix
$ws.Range("C2").Value = 4398.43020680935
$ws.Range("C3").Value = 4043.811468923269
$ws.Range("C4").Value = 4043.811468923269
$ws.Range("C5").Value = 3979.554928500677
$ws.Range("C6").Value = 3979.554928500677
$ws.Range("C7").Value = 3876.504347806812
$ws.Range("C8").Value = 3876.504347806812
$ws.Range("C9").Value = 3876.504347806812
$ws.Range("C10").Value = 3876.504347806812
$ws.Range("C11").Value = 3860.156550534577
$ws.Range("C12").Value = 3794.653020043877
